$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct ID column (G) typo: values were stored as plain numbers (1..7)
# but should be text ids "id01".."id07" (zero-padded, prefixed with "id").
for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G
    $num = [int]$cell.Value()
    $cell.Value = "id{0:D2}" -f $num
}

# Update the view selection to match the saved state
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H23").Select()
